$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coke Oven")
$ws.Range("D5").Formula = "=35*'Reference Values'!B18"
